$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header in D1 from "Ghi chú" to "Nhu cầu"
$ws.Range("D1").Value = "Nhu cầu"

# Update the active selection to D1 (was D2)
$ws.Range("D1").Select()
